$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $orig = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $orig
}

$ws.Range("D2").Value = "62.638.74"
$ws.Range("D3").Value = "2.439.39"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue $ws.Range("D5") "566.15"
$ws.Range("E5").Value = "  +0.47%  "
Set-TextValue $ws.Range("D6") "145.18"
$ws.Range("E6").Value = "  +1.75%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("E9").Value = "  +1.70%  "
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("E12").Value = "  +0.89%  "
Set-TextValue $ws.Range("D13") "26.86"
$ws.Range("E13").Value = "  +4.86%  "
$ws.Range("E14").Value = "  +4.97%  "
$ws.Range("D15").Value = "2.880.08"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "62.598.24"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "2.442.48"
$ws.Range("E17").Value = "  +1.25%  "
Set-TextValue $ws.Range("D18") "11.24"
$ws.Range("E18").Value = "  -0.30%  "
Set-TextValue $ws.Range("D19") "6.92"
$ws.Range("E19").Value = "  +0.49%  "
Set-TextValue $ws.Range("D20") "323.72"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("E21").Value = "  +0.35%  "
Set-TextValue $ws.Range("D22") "0.998"
$ws.Range("E22").Value = "  -0.21%  "
Set-TextValue $ws.Range("D23") "67.19"
$ws.Range("E23").Value = "  +1.68%  "
Set-TextValue $ws.Range("D24") "1.77"
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("E26").Value = "  +7.09%  "
Set-TextValue $ws.Range("D27") "567.72"
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("D28").Value = "2.553.62"
$ws.Range("E28").Value = "  +0.92%  "
Set-TextValue $ws.Range("D29") "0.999"
$ws.Range("E29").Value = "  -0.17%  "
Set-TextValue $ws.Range("D30") "8.37"
$ws.Range("E30").Value = "  +2.31%  "
$ws.Range("E31").Value = "  +2.55%  "
$ws.Range("E32").Value = "  -0.56%  "
Set-TextValue $ws.Range("D33") "1.87"
$ws.Range("E33").Value = "  -0.31%  "
Set-TextValue $ws.Range("D34") "1.55"
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("E35").Value = "  +3.42%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("E38").Value = "  -1.16%  "
Set-TextValue $ws.Range("D39") "18.79"
$ws.Range("E39").Value = "  +0.82%  "
Set-TextValue $ws.Range("D40") "148.65"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("E41").Value = "  +1.39%  "
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("E43").Value = "  +5.74%  "
Set-TextValue $ws.Range("D44") "148.13"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E45").Value = "  +1.30%  "
Set-TextValue $ws.Range("D46") "0.0535"
$ws.Range("E46").Value = "  +0.46%  "
Set-TextValue $ws.Range("D47") "20.44"
$ws.Range("E47").Value = "  +2.56%  "
Set-TextValue $ws.Range("D48") "0.602"
$ws.Range("E48").Value = "  +1.83%  "
Set-TextValue $ws.Range("D49") "0.0925"
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("E50").Value = "  +2.38%  "
$ws.Range("E51").Value = "  +2.03%  "
